$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.02258322285507441
$ws.Range("C2").Value = 0.3375848360084654
$ws.Range("D2").Value = 3.082599426703578
$ws.Range("E2").Value = 246.9852506941017
$ws.Range("G2").Value = 250.4280181796688

$ws.Range("B3").Value = 0.001754667048134761
$ws.Range("C3").Value = 0.05231270169004087
$ws.Range("D3").Value = 16.98373111632243
$ws.Range("E3").Value = 6.48142807727062
$ws.Range("G3").Value = 23.51922656233122
